$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 281, shifting existing rows 281:303 down to 282:304.
$ws.Rows.Item(281).Insert()

# New weekly record for Terminal La Palmera de La Serena - Ajo
$newRow = 281
$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44783
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112003
$ws.Cells.Item($newRow, 7).Value = "Ajo"
$ws.Cells.Item($newRow, 8).Value = "Chino"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 480
$ws.Cells.Item($newRow, 11).Value = 26000
$ws.Cells.Item($newRow, 12).Value = 26500
$ws.Cells.Item($newRow, 13).Value = 26250
$ws.Cells.Item($newRow, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item($newRow, 15).Value = "China"
$ws.Cells.Item($newRow, 16).Value = 2625
$ws.Cells.Item($newRow, 17).Value = 10
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
